# "se actualiza el diccionario de datos" - update the data dictionary column
# headers and drop the now-unused summary/balance columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the trailing "saldo" columns that are no longer part of the
#        dictionary (Saldo_Total_Debito, Saldo_Total_Credito, Saldo_Inicial).
#        Macro_campo_nivel_agregado (old column S) is kept and becomes the
#        new last column (Q) once Q:R are removed.
#        Delete from right to left so column letters don't shift underneath us.
$ws.Columns("T:T").Delete()   # Saldo_Inicial
$ws.Columns("Q:R").Delete()   # Saldo_Total_Debito, Saldo_Total_Credito

# --- 2. Rename the headers that make up the refreshed data dictionary.
$ws.Range("J1").Value = "Descripcion_Registro_Contable"
$ws.Range("L1").Value = "Identificacion_Tercero"
$ws.Range("N1").Value = "Descripcion_codigo_contable"

# --- 3. Adjust column widths to the refreshed dictionary's layout.
$ws.Columns("B:B").ColumnWidth = 15.72
$ws.Columns("C:C").ColumnWidth = 21.05
$ws.Columns("D:D").ColumnWidth = 14.05
$ws.Columns("E:E").ColumnWidth = 15.28
$ws.Columns("G:G").ColumnWidth = 22.83
$ws.Columns("H:H").ColumnWidth = 14.17
$ws.Columns("I:I").ColumnWidth = 23.17
$ws.Columns("J:J").ColumnWidth = 30.61
$ws.Columns("L:L").ColumnWidth = 21.61
$ws.Columns("M:M").ColumnWidth = 24.5
$ws.Columns("N:N").ColumnWidth = 27.5
$ws.Columns("Q:Q").ColumnWidth = 33.72

# --- 4. Update the selection to reflect where the author left off editing.
$ws.Range("Q6").Select()
